# Fix typo in the "Pod" bullet description: "A the smallest, ..." -> "The smallest, ..."
# The original run's text "A the smallest, schedulable resource that is managed by "
# is split into two runs:
#   1) "The "
#   2) "smallest, schedulable resource that is managed by "
# (the remaining runs for "kubelet" and " on the node" are left untouched)

$p = $ppt.ActivePresentation

$oldPhrase = "A the smallest, schedulable resource that is managed by "
$oldPrefix = "A the "
$newPrefix = "The "

$targetShape = $null
$targetSlide = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $sh = $sl.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $txt = $tf.TextRange.Text
                if ($txt.Contains($oldPhrase)) {
                    $targetShape = $sh
                    $targetSlide = $sl
                }
            }
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not find shape containing the target text"
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text
$charIdx = $fullText.IndexOf($oldPhrase)
if ($charIdx -lt 0) {
    throw "Target phrase not found"
}

# COM TextRange.Characters() is 1-based.
$startPos = $charIdx + 1

# Select just "A the " (the part being replaced) and retype it as "The ".
# Setting .Text on this sub-range splits the original run into two runs,
# leaving the untouched remainder ("smallest, schedulable resource that is
# managed by ") as its own run with the original formatting/dirty state.
$prefixRange = $tr.Characters($startPos, $oldPrefix.Length)
$prefixRange.Text = $newPrefix
